$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet logs timestamped ranking samples. A new sample taken on
# 2026/02/15 (Sunday) at hour 16 with ranking 201 needs to be inserted
# in date order, right after the existing 2026/02/15 13:00 entry and
# before the 2026/12/29 entries that currently start at row 825.
#
# Insert a new row at row 825; every following row (825-866) shifts
# down by one (826-867), and the sheet's used range grows to A1:D867.
$ws.Rows.Item(825).Insert()

# Column A stores the date as literal text (e.g. "2026/02/15"), not a
# real Excel date serial. Typing a date-shaped string through COM
# auto-converts it to a date, so force Text formatting first, write the
# value, then clear the formatting again so the cell ends up with the
# same (default/no) style as every other data row in the sheet.
$ws.Cells.Item(825, 1).NumberFormat = "@"
$ws.Cells.Item(825, 1).Value = "2026/02/15"
$ws.Cells.Item(825, 1).ClearFormats()

$ws.Cells.Item(825, 2).Value = "日"
$ws.Cells.Item(825, 3).Value = 16
$ws.Cells.Item(825, 4).Value = 201
